# Auto-generated edit script for cryptos.xlsx update
# Applies the cell-level changes described by the commit diff:
#   "Updated cryptos list on Thu Apr 18 10:32:20 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells receive values that look numeric ("542.24", "0.999", ...).
# The source sheet stores every Price/Volume figure as literal text, so
# force a Text format before assigning to stop Excel from auto-converting
# them into numbers.
$textCells = @("D5","D6","D9","D10","D11","D12","D13","D14","D19","D20","D21","D22","D23","D24","D25","D26","D28","D29","D31","D33","D34","D35","D36","D38","D39","D42","D43","D44","D45","D46","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, in sheet order.
$ws.Range("D2").Value = "61.683.42"
$ws.Range("E2").Value = "  -2.04%  "
$ws.Range("D3").Value = "3.005.48"
$ws.Range("E3").Value = "  -2.10%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "542.24"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").Value = "131.34"
$ws.Range("E6").Value = "  -4.69%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "2.998.52"
$ws.Range("E8").Value = "  -2.11%  "
$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").Value = "6.07"
$ws.Range("E10").Value = "  -3.16%  "
$ws.Range("D11").Value = "0.145"
$ws.Range("E11").Value = "  -7.32%  "
$ws.Range("D12").Value = "0.445"
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("D13").Value = "34.28"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").Value = "0.0000219"
$ws.Range("E14").Value = "  -2.03%  "
$ws.Range("D15").Value = "3.489.54"
$ws.Range("E15").Value = "  -2.17%  "
$ws.Range("D16").Value = "61.649.64"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("E17").Value = "  -2.89%  "
$ws.Range("D18").Value = "3.002.01"
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("D19").Value = "6.60"
$ws.Range("D20").Value = "483.21"
$ws.Range("E20").Value = "  +2.64%  "
$ws.Range("D21").Value = "13.21"
$ws.Range("E21").Value = "  -2.66%  "
$ws.Range("D22").Value = "0.666"
$ws.Range("E22").Value = "  -4.53%  "
$ws.Range("D23").Value = "6.93"
$ws.Range("E23").Value = "  -1.71%  "
$ws.Range("D24").Value = "82.16"
$ws.Range("E24").Value = "  +4.57%  "
$ws.Range("D25").Value = "11.91"
$ws.Range("E25").Value = "  -1.83%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D28").Value = "7.64"
$ws.Range("E28").Value = "  -3.07%  "
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("D31").Value = "25.64"
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("E32").Value = "  -3.46%  "
$ws.Range("D33").Value = "5.59"
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").Value = "2.32"
$ws.Range("E34").Value = "  +0.77%  "
$ws.Range("D35").Value = "54.78"
$ws.Range("E35").Value = "  -6.91%  "
$ws.Range("D36").Value = "5.83"
$ws.Range("E36").Value = "  -2.63%  "
$ws.Range("D37").Value = "3.133.50"
$ws.Range("E37").Value = "  -3.68%  "
$ws.Range("D38").Value = "435.73"
$ws.Range("E38").Value = "  -10.37%  "
$ws.Range("D39").Value = "0.0792"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("E40").Value = "  -4.76%  "
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").Value = "8.05"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").Value = "2.40"
$ws.Range("E43").Value = "  -6.76%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "26.21"
$ws.Range("E44").Value = "  +3.63%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").Value = "0.240"
$ws.Range("E46").Value = "  -4.59%  "
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").Value = "1.93"
$ws.Range("E48").Value = "  -3.67%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "115.15"
$ws.Range("E49").Value = "  -6.25%  "
$ws.Range("D50").Value = "1.28"
$ws.Range("E50").Value = "  +4.14%  "
$ws.Range("D51").Value = "0.0₃0483"
$ws.Range("E51").Value = "  -7.63%  "

Write-Host "Applied cryptos list update"
